$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AF ("Business TMP") and AG ("CATEGORIA") hold sample wildcard-matching
# test data. Replace the old "abcde / abcd qualsiasi cosa / mikecccc / ..."
# placeholder values with the new "Caso D / Caso 1 / Pippo N / Pluto N" set,
# and reshuffle the existing wildcard patterns across rows 3-24.
$values = @{
    3  = @("*cccc",   "abcd*")
    4  = @("*ilmn",   "ef*gh")
    5  = @("aaaa*",   "Caso D")
    6  = @("abcd*",   "Caso 1")
    7  = @("bb*bb",   "bb*bb")
    8  = @("Caso 1",  "aaaa*")
    9  = @("Caso D",  "*ilmn")
    10 = @("ef*gh",   "*cccc")
    11 = @("*cccc",   "abcd*")
    12 = @("*ilmn",   "ef*gh")
    13 = @("aaaa*",   "Caso D")
    14 = @("abcd*",   "Caso 1")
    15 = @("bb*bb",   "bb*bb")
    16 = @("Caso 1",  "aaaa*")
    17 = @("Caso D",  "*ilmn")
    18 = @("ef*gh",   "*cccc")
    19 = @("Pluto 1", "Pippo 1")
    20 = @("Pluto 2", "Pippo 2")
    21 = @("Pluto 3", "Pippo 3")
    22 = @("Pippo 1", "Pluto 1")
    23 = @("Pippo 2", "Pluto 2")
    24 = @("Pippo 3", "Pluto 3")
}

# Write column AG first, then column AF (both top-to-bottom), matching the
# order the new sample values were authored in.
foreach ($row in ($values.Keys | Sort-Object)) {
    $ws.Range("AG$row").Value = $values[$row][1]
}
foreach ($row in ($values.Keys | Sort-Object)) {
    $ws.Range("AF$row").Value = $values[$row][0]
}

# Update view state: scroll so column M is the left-most visible column, and
# move the active selection to a single cell.
$ws.Activate()
$ws.Range("AG23").Select()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
